# Custom Report Feature Updated through IReporter
# Updates the recorded Start Time / End Time timestamps for the test
# cases on the active ResultSummary sheet, reflecting a new IReporter run.
# testCase1 and testCase2 shared the same "04-Jan-2020 04:48:19" start/end
# timestamp, so both move to the new run's single timestamp together.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# testCase1 (row 3) and testCase2 (row 4): Start Time (B) and End Time (C)
# all become the new run's single timestamp.
$ws.Range("B3").Value = "05-Jan-2020 02:07:31"
$ws.Range("C3").Value = "05-Jan-2020 02:07:31"
$ws.Range("B4").Value = "05-Jan-2020 02:07:31"
$ws.Range("C4").Value = "05-Jan-2020 02:07:31"

# test01 (row 6): Start Time (B6) and End Time (C6) both become the new
# run's single timestamp.
$ws.Range("B6").Value = "05-Jan-2020 02:07:40"
$ws.Range("C6").Value = "05-Jan-2020 02:07:40"
